$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 10 ("MIXED" restriction) -----------------------------
# Shift the existing last row (row 11: NULL / additionalProperties / Must NULL)
# down to row 12, leaving row 11 empty (a gap row) just like the target layout.
$ws.Range("B11:E11").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Bring in the formatting that the new row needs from cells that already carry it,
# so the existing font/fill/border/alignment records in the workbook get reused
# instead of ad-hoc ones being fabricated.

# C10 ("RESTRICTION" column) reuses the wrapped/bordered text style used by C9.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# D10/E10 reuse the plain bordered style used throughout column D/E (e.g. D3:E3).
$ws.Range("D3:E3").Copy()
$ws.Range("D10:E10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# B10 ("TYPE" column) starts from the bordered style used by B3, then gets a bold
# Arial Narrow 10pt font to stand out as a section/merged-type heading ("MIXED").
$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B10").Font.Bold = $true
$ws.Range("B10").Font.Name = "Arial Narrow"
$ws.Range("B10").Font.Size = 10

$ws.Application.CutCopyMode = $false

# Fill in the new row's content.
$ws.Range("C10").Value = "`"type`": [ `"integer`", `"string`" ]," + [char]10 + "`"minimum`": 1," + [char]10 + "`"minLength`": 1"
$ws.Range("B10").Value = "MIXED"

# Row height matches the other 3-line wrapped restriction rows.
$ws.Rows("10:10").RowHeight = 38.25

# --- Update the sheet view selection ---------------------------------------
# The bottom-right (scrollable) pane's active cell moves to B11 (the blank
# gap row right below the newly inserted MIXED row).
$ws.Range("B11").Select()
